$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.109201
$ws.Range("H2").Value = 0.327603
$ws.Range("I2").Value = 0.07562717345335074
$ws.Range("J2").Value = 0.07562717345335074
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1311436666666667
$ws.Range("N2").Value = 0.393431
$ws.Range("O2").Value = 0.02663441993971509
$ws.Range("P2").Value = 0.02663441993971509
$ws.Range("Q2").Value = 0.01432101954366666
$ws.Range("R2").Value = 0.128889175893
$ws.Range("S2").Value = 0.002014285896610216
$ws.Range("T2").Value = 0.002014285896610216
$ws.Range("G3").Value = 0.109201
$ws.Range("H3").Value = 0.327603
$ws.Range("I3").Value = 0.07562717345335074
$ws.Range("J3").Value = 0.07562717345335074
$ws.Range("O3").Value = 0.06149297381279183
$ws.Range("P3").Value = 0.06149297381279183
$ws.Range("Q3").Value = 0.03306406078166667
$ws.Range("R3").Value = 0.297576547035
$ws.Range("S3").Value = 0.004650539796702362
$ws.Range("T3").Value = 0.004650539796702362
$ws.Range("G4").Value = 0.109201
$ws.Range("H4").Value = 0.327603
$ws.Range("I4").Value = 0.07562717345335074
$ws.Range("J4").Value = 0.07562717345335074
$ws.Range("M4").Value = 4.009307333333333
$ws.Range("N4").Value = 12.027922
$ws.Range("O4").Value = 0.8142640654908683
$ws.Range("P4").Value = 0.8142640654908684
$ws.Range("Q4").Value = 0.4378203701073333
$ws.Range("R4").Value = 3.940383330966
$ws.Range("S4").Value = 0.06158048971770844
$ws.Range("T4").Value = 0.06158048971770845
$ws.Range("G5").Value = 0.109201
$ws.Range("H5").Value = 0.327603
$ws.Range("I5").Value = 0.07562717345335074
$ws.Range("J5").Value = 0.07562717345335074
$ws.Range("M5").Value = 0.480609
$ws.Range("N5").Value = 1.441827
$ws.Range("O5").Value = 0.09760854075662465
$ws.Range("P5").Value = 0.09760854075662465
$ws.Range("Q5").Value = 0.052482983409
$ws.Range("R5").Value = 0.472346850681
$ws.Range("S5").Value = 0.007381858042329707
$ws.Range("T5").Value = 0.007381858042329707
$ws.Range("G6").Value = 0.7328223333333334
$ws.Range("I6").Value = 0.5075162472274908
$ws.Range("J6").Value = 0.5075162472274908
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1311436666666667
$ws.Range("N6").Value = 0.393431
$ws.Range("O6").Value = 0.02663441993971509
$ws.Range("P6").Value = 0.02663441993971509
$ws.Range("Q6").Value = 0.09610500780855555
$ws.Range("R6").Value = 0.864945070277
$ws.Range("S6").Value = 0.01351740085488525
$ws.Range("T6").Value = 0.01351740085488525
$ws.Range("G7").Value = 0.7328223333333334
$ws.Range("I7").Value = 0.5075162472274908
$ws.Range("J7").Value = 0.5075162472274908
$ws.Range("O7").Value = 0.06149297381279183
$ws.Range("P7").Value = 0.06149297381279183
$ws.Range("Q7").Value = 0.2218851674572222
$ws.Range("S7").Value = 0.03120868330032648
$ws.Range("T7").Value = 0.03120868330032648
$ws.Range("G8").Value = 0.7328223333333334
$ws.Range("I8").Value = 0.5075162472274908
$ws.Range("J8").Value = 0.5075162472274908
$ws.Range("M8").Value = 4.009307333333333
$ws.Range("N8").Value = 12.027922
$ws.Range("O8").Value = 0.8142640654908683
$ws.Range("P8").Value = 0.8142640654908684
$ws.Range("Q8").Value = 2.938109955063778
$ws.Range("R8").Value = 26.442989595574
$ws.Range("S8").Value = 0.4132522427701252
$ws.Range("T8").Value = 0.4132522427701253
$ws.Range("G9").Value = 0.7328223333333334
$ws.Range("I9").Value = 0.5075162472274908
$ws.Range("J9").Value = 0.5075162472274908
$ws.Range("M9").Value = 0.480609
$ws.Range("N9").Value = 1.441827
$ws.Range("O9").Value = 0.09760854075662465
$ws.Range("P9").Value = 0.09760854075662465
$ws.Range("Q9").Value = 0.352201008801
$ws.Range("R9").Value = 3.169809079209
$ws.Range("S9").Value = 0.04953792030215372
$ws.Range("T9").Value = 0.04953792030215372
$ws.Range("G10").Value = 0.498848
$ws.Range("H10").Value = 1.496544
$ws.Range("I10").Value = 0.3454772778899196
$ws.Range("J10").Value = 0.3454772778899197
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1311436666666667
$ws.Range("N10").Value = 0.393431
$ws.Range("O10").Value = 0.02663441993971509
$ws.Range("P10").Value = 0.02663441993971509
$ws.Range("Q10").Value = 0.06542075582933334
$ws.Range("R10").Value = 0.5887868024639999
$ws.Range("S10").Value = 0.009201586898949767
$ws.Range("T10").Value = 0.009201586898949768
$ws.Range("G11").Value = 0.498848
$ws.Range("H11").Value = 1.496544
$ws.Range("I11").Value = 0.3454772778899196
$ws.Range("J11").Value = 0.3454772778899197
$ws.Range("O11").Value = 0.06149297381279183
$ws.Range("P11").Value = 0.06149297381279183
$ws.Range("Q11").Value = 0.1510420288533333
$ws.Range("R11").Value = 1.35937825968
$ws.Range("S11").Value = 0.02124442520219944
$ws.Range("T11").Value = 0.02124442520219944
$ws.Range("G12").Value = 0.498848
$ws.Range("H12").Value = 1.496544
$ws.Range("I12").Value = 0.3454772778899196
$ws.Range("J12").Value = 0.3454772778899197
$ws.Range("M12").Value = 4.009307333333333
$ws.Range("N12").Value = 12.027922
$ws.Range("O12").Value = 0.8142640654908683
$ws.Range("P12").Value = 0.8142640654908684
$ws.Range("Q12").Value = 2.000034944618667
$ws.Range("R12").Value = 18.000314501568
$ws.Range("S12").Value = 0.2813097328293644
$ws.Range("T12").Value = 0.2813097328293645
$ws.Range("G13").Value = 0.498848
$ws.Range("H13").Value = 1.496544
$ws.Range("I13").Value = 0.3454772778899196
$ws.Range("J13").Value = 0.3454772778899197
$ws.Range("M13").Value = 0.480609
$ws.Range("N13").Value = 1.441827
$ws.Range("O13").Value = 0.09760854075662465
$ws.Range("P13").Value = 0.09760854075662465
$ws.Range("Q13").Value = 0.239750838432
$ws.Range("R13").Value = 2.157757545888
$ws.Range("S13").Value = 0.03372153295940596
$ws.Range("T13").Value = 0.03372153295940596
$ws.Range("G14").Value = 0.1030673333333333
$ws.Range("H14").Value = 0.309202
$ws.Range("I14").Value = 0.07137930142923891
$ws.Range("J14").Value = 0.07137930142923891
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1311436666666667
$ws.Range("N14").Value = 0.393431
$ws.Range("O14").Value = 0.02663441993971509
$ws.Range("P14").Value = 0.02663441993971509
$ws.Range("Q14").Value = 0.01351662800688889
$ws.Range("R14").Value = 0.121649652062
$ws.Range("S14").Value = 0.001901146289269855
$ws.Range("T14").Value = 0.001901146289269855
$ws.Range("G15").Value = 0.1030673333333333
$ws.Range("H15").Value = 0.309202
$ws.Range("I15").Value = 0.07137930142923891
$ws.Range("J15").Value = 0.07137930142923891
$ws.Range("O15").Value = 0.06149297381279183
$ws.Range("P15").Value = 0.06149297381279183
$ws.Range("Q15").Value = 0.03120689896555556
$ws.Range("R15").Value = 0.28086209069
$ws.Range("S15").Value = 0.004389325513563564
$ws.Range("T15").Value = 0.004389325513563564
$ws.Range("G16").Value = 0.1030673333333333
$ws.Range("H16").Value = 0.309202
$ws.Range("I16").Value = 0.07137930142923891
$ws.Range("J16").Value = 0.07137930142923891
$ws.Range("M16").Value = 4.009307333333333
$ws.Range("N16").Value = 12.027922
$ws.Range("O16").Value = 0.8142640654908683
$ws.Range("P16").Value = 0.8142640654908684
$ws.Range("Q16").Value = 0.4132286153604444
$ws.Range("R16").Value = 3.719057538244
$ws.Range("S16").Value = 0.05812160017367022
$ws.Range("T16").Value = 0.05812160017367023
$ws.Range("G17").Value = 0.1030673333333333
$ws.Range("H17").Value = 0.309202
$ws.Range("I17").Value = 0.07137930142923891
$ws.Range("J17").Value = 0.07137930142923891
$ws.Range("M17").Value = 0.480609
$ws.Range("N17").Value = 1.441827
$ws.Range("O17").Value = 0.09760854075662465
$ws.Range("P17").Value = 0.09760854075662465
$ws.Range("Q17").Value = 0.049535088006
$ws.Range("R17").Value = 0.445815792054
$ws.Range("S17").Value = 0.006967229452735263
$ws.Range("T17").Value = 0.006967229452735263
